{"js": "// 1) Update the MSc thesis topic title text.\nconst oldTitle =\n  \"Usability and scalability of textual notations vs graphical notations for model transformation languages\";\nconst newTitle =\n  \"Design and implementation of a textual notation for a graphical model transformation language\";\n\nconst titleHits = context.document.body.search(oldTitle, { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length === 0) {\n  throw new Error(\"Could not find the thesis title text to replace.\");\n}\ntitleHits.items[0].insertText(newTitle, \"Replace\");\nawait context.sync();\n\n// 2) Remove the empty paragraph + the \"In a broader context...\" paragraph\n//    that followed the \"...front-end used to develop the model\n//    transformation specification itself.\" paragraph.\nconst anchorHits = context.document.body.search(\n  \"front-end used to develop the model transformation specification itself.\",\n  { matchCase: true }\n);\nanchorHits.load(\"items\");\nawait context.sync();\n\nif (anchorHits.items.length === 0) {\n  throw new Error(\"Could not find the anchor paragraph for the deletion.\");\n}\n\nconst anchorParagraph = anchorHits.items[0].paragraphs.getFirst();\nconst emptyParagraph = anchorParagraph.getNext();\nconst contextParagraph = emptyParagraph.getNext();\n\nemptyParagraph.load(\"text\");\ncontextParagraph.load(\"text\");\nawait context.sync();\n\nif (!contextParagraph.text.startsWith(\"In a broader context\")) {\n  throw new Error(\"Unexpected document structure near the deletion target.\");\n}\n\ncontextParagraph.delete();\nemptyParagraph.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the MSc thesis topic title text.\n$oldTitle = \"Usability and scalability of textual notations vs graphical notations for model transformation languages\"\n$newTitle = \"Design and implementation of a textual notation for a graphical model transformation language\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute($oldTitle, $false, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)\n\n# 2) Remove the empty paragraph + the \"In a broader context...\" paragraph\n#    that followed the \"...front-end used to develop the model\n#    transformation specification itself.\" paragraph.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*front-end used to develop the model transformation specification itself.*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the anchor paragraph for the deletion.\"\n}\n\n$emptyParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n$contextParagraph = $d.Paragraphs.Item($anchorIndex + 2)\n\nif ($contextParagraph.Range.Text.StartsWith(\"In a broader context\")) {\n    $contextParagraph.Range.Delete()\n    $emptyParagraph.Range.Delete()\n} else {\n    throw \"Unexpected document structure near the deletion target.\"\n}\n"}
